# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-27 23:14:16
#
# Updates the "Recorded By" column (G) on the "Session Analysis Results"
# sheet: the list of recorders in several rows is re-ordered (the values
# themselves are unchanged, only their order within the comma-separated
# list is swapped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the 2-item "Recorded By" list (G column) needs its two
# entries swapped, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# or "admin@admin.com, dnasr281@gmail.com" -> "dnasr281@gmail.com, admin@admin.com"
$rowsSwap2 = @(
    3, 6, 10, 12, 13, 14, 15, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 38, 39, 40,
    41, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 64, 65, 66, 67, 70, 71, 72, 73, 74,
    76, 78, 83, 84, 85, 86, 87, 90, 92, 99, 101, 109, 110, 111, 112, 113, 116, 118,
    125, 127, 135, 136, 137, 138, 139, 142, 144, 151, 153
)

# Rows where the 3-item "Recorded By" list has its last two entries
# swapped, e.g. "backup@backdoor.com, System, system" -> "backup@backdoor.com, system, System"
$rowsSwapLast2Of3 = @(2, 28, 54)

$col = 7  # column G = "Recorded By"

foreach ($r in $rowsSwap2) {
    $cell = $ws.Cells.Item($r, $col)
    $parts = $cell.Value2 -split ", "
    $cell.Value = $parts[1] + ", " + $parts[0]
}

foreach ($r in $rowsSwapLast2Of3) {
    $cell = $ws.Cells.Item($r, $col)
    $parts = $cell.Value2 -split ", "
    $cell.Value = $parts[0] + ", " + $parts[2] + ", " + $parts[1]
}
